# Regenerate orders with updated distance/size codes.
#
# The experiment's distance codes and one size code were renumbered:
#   D80 -> D86
#   D51 -> D55
#   D64 -> D69
#   S30 -> S31
#
# These tokens appear as whole values (e.g. "D80") and embedded inside
# composite labels and filenames (e.g. "Face13_D80_S25",
# "Face13_D80_S25_l.png"). A substring Replace across the used range
# updates every occurrence consistently.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ur = $ws.UsedRange

# Order matters only in that none of the replacement targets collide with
# another rule's source text, so a simple sequential pass is safe.
$ur.Replace("D80", "D86") | Out-Null
$ur.Replace("D51", "D55") | Out-Null
$ur.Replace("D64", "D69") | Out-Null
$ur.Replace("S30", "S31") | Out-Null
